{"js": "// Add a \"Recommendation Type\" / \"Lighting\" row to the top of the\n// \"Summary of Estimated Savings and Implementation Costs\" table, and\n// restore the table's (slightly wider) column widths.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Insert the new row as the first row of the table.\ntable.addRows(\"Start\", 1, [[\"Recommendation Type\", \"Lighting\"]]);\nawait context.sync();\n\n// Setting columnWidth on a cell resizes every cell in that column (and the\n// <w:tblGrid> entry), so two writes re-establish the new widths for the\n// whole table: 3703 dxa (185.15 pt) for column 1, 2222 dxa (111.1 pt) for\n// column 2 (dxa values are twentieths of a point).\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst firstRowCells = rows.items[0].cells;\nfirstRowCells.load(\"items\");\nawait context.sync();\n\nfirstRowCells.items[0].columnWidth = 3703 / 20;\nfirstRowCells.items[1].columnWidth = 2222 / 20;\nawait context.sync();\n", "ps1": "# Add a \"Recommendation Type\" / \"Lighting\" row to the top of the\n# \"Summary of Estimated Savings and Implementation Costs\" table, and\n# restore the table's (slightly wider) column widths.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Insert a new row before the current first row.\n$firstRow = $t.Rows.Item(1)\n$newRow = $t.Rows.Add($firstRow)\n$newRow.Cells.Item(1).Range.Text = \"Recommendation Type\"\n$newRow.Cells.Item(2).Range.Text = \"Lighting\"\n\n# Re-apply the table's column widths (3703 dxa / 2222 dxa, i.e. 185.15pt /\n# 111.1pt since a point is 20 dxa) to every row, including the new one.\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    $row = $t.Rows.Item($i)\n    $row.Cells.Item(1).Width = 185.15\n    $row.Cells.Item(2).Width = 111.1\n}\n"}
